# Natmi following Dr Hou advice
# Add "ECs" as a sending-cluster category (mirrors the already-present
# "ECs" target-cluster rows), expanding the 3x4 grid to a 4x4 grid.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Hgf"
$ws.Cells.Item(2, 3).Value = "Cd44"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 8.265822
$ws.Cells.Item(2, 8).Value = 24.797466
$ws.Cells.Item(2, 9).Value = 0.2082338764513023
$ws.Cells.Item(2, 10).Value = 0.2082338764513023
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 239.0839323333333
$ws.Cells.Item(2, 14).Value = 717.251797
$ws.Cells.Item(2, 15).Value = 0.4086975387666237
$ws.Cells.Item(2, 16).Value = 0.4086975387666237
$ws.Cells.Item(2, 17).Value = 1976.225227727378
$ws.Cells.Item(2, 18).Value = 17786.0270495464
$ws.Cells.Item(2, 19).Value = 0.08510467279348043
$ws.Cells.Item(2, 20).Value = 0.08510467279348043

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Hgf"
$ws.Cells.Item(3, 3).Value = "Cd44"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 8.265822
$ws.Cells.Item(3, 8).Value = 24.797466
$ws.Cells.Item(3, 9).Value = 0.2082338764513023
$ws.Cells.Item(3, 10).Value = 0.2082338764513023
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 117.0512696666667
$ws.Cells.Item(3, 14).Value = 351.153809
$ws.Cells.Item(3, 15).Value = 0.2000910950200451
$ws.Cells.Item(3, 16).Value = 0.2000910950200451
$ws.Cells.Item(3, 17).Value = 967.524959938666
$ws.Cells.Item(3, 18).Value = 8707.724639447995
$ws.Cells.Item(3, 19).Value = 0.04166574435940985
$ws.Cells.Item(3, 20).Value = 0.04166574435940985

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Hgf"
$ws.Cells.Item(4, 3).Value = "Cd44"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 8.265822
$ws.Cells.Item(4, 8).Value = 24.797466
$ws.Cells.Item(4, 9).Value = 0.2082338764513023
$ws.Cells.Item(4, 10).Value = 0.2082338764513023
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 171.15883
$ws.Cells.Item(4, 14).Value = 513.47649
$ws.Cells.Item(4, 15).Value = 0.2925842480357353
$ws.Cells.Item(4, 16).Value = 0.2925842480357353
$ws.Cells.Item(4, 17).Value = 1414.76842250826
$ws.Cells.Item(4, 18).Value = 12732.91580257434
$ws.Cells.Item(4, 19).Value = 0.06092595215707048
$ws.Cells.Item(4, 20).Value = 0.06092595215707048

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Hgf"
$ws.Cells.Item(5, 3).Value = "Cd44"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 8.265822
$ws.Cells.Item(5, 8).Value = 24.797466
$ws.Cells.Item(5, 9).Value = 0.2082338764513023
$ws.Cells.Item(5, 10).Value = 0.2082338764513023
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 57.695868
$ws.Cells.Item(5, 14).Value = 173.087604
$ws.Cells.Item(5, 15).Value = 0.09862711817759588
$ws.Cells.Item(5, 16).Value = 0.09862711817759588
$ws.Cells.Item(5, 17).Value = 476.903775023496
$ws.Cells.Item(5, 18).Value = 4292.133975211464
$ws.Cells.Item(5, 19).Value = 0.02053750714134149
$ws.Cells.Item(5, 20).Value = 0.02053750714134149

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Hgf"
$ws.Cells.Item(6, 3).Value = "Cd44"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 11.28595333333333
$ws.Cells.Item(6, 8).Value = 33.85786
$ws.Cells.Item(6, 9).Value = 0.2843174958338682
$ws.Cells.Item(6, 10).Value = 0.2843174958338682
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 239.0839323333333
$ws.Cells.Item(6, 14).Value = 717.251797
$ws.Cells.Item(6, 15).Value = 0.4086975387666237
$ws.Cells.Item(6, 16).Value = 0.4086975387666237
$ws.Cells.Item(6, 17).Value = 2698.290103063825
$ws.Cells.Item(6, 18).Value = 24284.61092757442
$ws.Cells.Item(6, 19).Value = 0.1161998607755917
$ws.Cells.Item(6, 20).Value = 0.1161998607755917

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Hgf"
$ws.Cells.Item(7, 3).Value = "Cd44"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 11.28595333333333
$ws.Cells.Item(7, 8).Value = 33.85786
$ws.Cells.Item(7, 9).Value = 0.2843174958338682
$ws.Cells.Item(7, 10).Value = 0.2843174958338682
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 117.0512696666667
$ws.Cells.Item(7, 14).Value = 351.153809
$ws.Cells.Item(7, 15).Value = 0.2000910950200451
$ws.Cells.Item(7, 16).Value = 0.2000910950200451
$ws.Cells.Item(7, 17).Value = 1321.035167065416
$ws.Cells.Item(7, 18).Value = 11889.31650358874
$ws.Cells.Item(7, 19).Value = 0.0568893990747558
$ws.Cells.Item(7, 20).Value = 0.0568893990747558

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Hgf"
$ws.Cells.Item(8, 3).Value = "Cd44"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 11.28595333333333
$ws.Cells.Item(8, 8).Value = 33.85786
$ws.Cells.Item(8, 9).Value = 0.2843174958338682
$ws.Cells.Item(8, 10).Value = 0.2843174958338682
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 171.15883
$ws.Cells.Item(8, 14).Value = 513.47649
$ws.Cells.Item(8, 15).Value = 0.2925842480357353
$ws.Cells.Item(8, 16).Value = 0.2925842480357353
$ws.Cells.Item(8, 17).Value = 1931.690567967933
$ws.Cells.Item(8, 18).Value = 17385.2151117114
$ws.Cells.Item(8, 19).Value = 0.08318682072195564
$ws.Cells.Item(8, 20).Value = 0.08318682072195564

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Hgf"
$ws.Cells.Item(9, 3).Value = "Cd44"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 11.28595333333333
$ws.Cells.Item(9, 8).Value = 33.85786
$ws.Cells.Item(9, 9).Value = 0.2843174958338682
$ws.Cells.Item(9, 10).Value = 0.2843174958338682
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 57.695868
$ws.Cells.Item(9, 14).Value = 173.087604
$ws.Cells.Item(9, 15).Value = 0.09862711817759588
$ws.Cells.Item(9, 16).Value = 0.09862711817759588
$ws.Cells.Item(9, 17).Value = 651.15287377416
$ws.Cells.Item(9, 18).Value = 5860.375863967441
$ws.Cells.Item(9, 19).Value = 0.02804141526156505
$ws.Cells.Item(9, 20).Value = 0.02804141526156505

# Row 10
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Hgf"
$ws.Cells.Item(10, 3).Value = "Cd44"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 19.51551966666667
$ws.Cells.Item(10, 8).Value = 58.546559
$ws.Cells.Item(10, 9).Value = 0.4916380138783083
$ws.Cells.Item(10, 10).Value = 0.4916380138783083
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 239.0839323333333
$ws.Cells.Item(10, 14).Value = 717.251797
$ws.Cells.Item(10, 15).Value = 0.4086975387666237
$ws.Cells.Item(10, 16).Value = 0.4086975387666237
$ws.Cells.Item(10, 17).Value = 4665.847183435169
$ws.Cells.Item(10, 18).Value = 41992.62465091653
$ws.Cells.Item(10, 19).Value = 0.2009312462361758
$ws.Cells.Item(10, 20).Value = 0.2009312462361758

# Row 11
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Hgf"
$ws.Cells.Item(11, 3).Value = "Cd44"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 19.51551966666667
$ws.Cells.Item(11, 8).Value = 58.546559
$ws.Cells.Item(11, 9).Value = 0.4916380138783083
$ws.Cells.Item(11, 10).Value = 0.4916380138783083
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 117.0512696666667
$ws.Cells.Item(11, 14).Value = 351.153809
$ws.Cells.Item(11, 15).Value = 0.2000910950200451
$ws.Cells.Item(11, 16).Value = 0.2000910950200451
$ws.Cells.Item(11, 17).Value = 2284.316355188137
$ws.Cells.Item(11, 18).Value = 20558.84719669323
$ws.Cells.Item(11, 19).Value = 0.09837238855039082
$ws.Cells.Item(11, 20).Value = 0.09837238855039082

# Row 12
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Hgf"
$ws.Cells.Item(12, 3).Value = "Cd44"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 19.51551966666667
$ws.Cells.Item(12, 8).Value = 58.546559
$ws.Cells.Item(12, 9).Value = 0.4916380138783083
$ws.Cells.Item(12, 10).Value = 0.4916380138783083
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 171.15883
$ws.Cells.Item(12, 14).Value = 513.47649
$ws.Cells.Item(12, 15).Value = 0.2925842480357353
$ws.Cells.Item(12, 16).Value = 0.2925842480357353
$ws.Cells.Item(12, 17).Value = 3340.253512988656
$ws.Cells.Item(12, 18).Value = 30062.28161689791
$ws.Cells.Item(12, 19).Value = 0.1438455385963672
$ws.Cells.Item(12, 20).Value = 0.1438455385963672

# Row 13
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Hgf"
$ws.Cells.Item(13, 3).Value = "Cd44"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 19.51551966666667
$ws.Cells.Item(13, 8).Value = 58.546559
$ws.Cells.Item(13, 9).Value = 0.4916380138783083
$ws.Cells.Item(13, 10).Value = 0.4916380138783083
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 57.695868
$ws.Cells.Item(13, 14).Value = 173.087604
$ws.Cells.Item(13, 15).Value = 0.09862711817759588
$ws.Cells.Item(13, 16).Value = 0.09862711817759588
$ws.Cells.Item(13, 17).Value = 1125.964846639404
$ws.Cells.Item(13, 18).Value = 10133.68361975464
$ws.Cells.Item(13, 19).Value = 0.04848884049537443
$ws.Cells.Item(13, 20).Value = 0.04848884049537443

# Row 14
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Hgf"
$ws.Cells.Item(14, 3).Value = "Cd44"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0.6276006666666666
$ws.Cells.Item(14, 8).Value = 1.882802
$ws.Cells.Item(14, 9).Value = 0.01581061383652123
$ws.Cells.Item(14, 10).Value = 0.01581061383652123
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 239.0839323333333
$ws.Cells.Item(14, 14).Value = 717.251797
$ws.Cells.Item(14, 15).Value = 0.4086975387666237
$ws.Cells.Item(14, 16).Value = 0.4086975387666237
$ws.Cells.Item(14, 17).Value = 150.0492353216882
$ws.Cells.Item(14, 18).Value = 1350.443117895194
$ws.Cells.Item(14, 19).Value = 0.006461758961375751
$ws.Cells.Item(14, 20).Value = 0.006461758961375751

# Row 15
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Hgf"
$ws.Cells.Item(15, 3).Value = "Cd44"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.6276006666666666
$ws.Cells.Item(15, 8).Value = 1.882802
$ws.Cells.Item(15, 9).Value = 0.01581061383652123
$ws.Cells.Item(15, 10).Value = 0.01581061383652123
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 117.0512696666667
$ws.Cells.Item(15, 14).Value = 351.153809
$ws.Cells.Item(15, 15).Value = 0.2000910950200451
$ws.Cells.Item(15, 16).Value = 0.2000910950200451
$ws.Cells.Item(15, 17).Value = 73.46145487697977
$ws.Cells.Item(15, 18).Value = 661.153093892818
$ws.Cells.Item(15, 19).Value = 0.003163563035488608
$ws.Cells.Item(15, 20).Value = 0.003163563035488608

# Row 16
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Hgf"
$ws.Cells.Item(16, 3).Value = "Cd44"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0.6276006666666666
$ws.Cells.Item(16, 8).Value = 1.882802
$ws.Cells.Item(16, 9).Value = 0.01581061383652123
$ws.Cells.Item(16, 10).Value = 0.01581061383652123
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 171.15883
$ws.Cells.Item(16, 14).Value = 513.47649
$ws.Cells.Item(16, 15).Value = 0.2925842480357353
$ws.Cells.Item(16, 16).Value = 0.2925842480357353
$ws.Cells.Item(16, 17).Value = 107.4193958138867
$ws.Cells.Item(16, 18).Value = 966.7745623249799
$ws.Cells.Item(16, 19).Value = 0.004625936560341956
$ws.Cells.Item(16, 20).Value = 0.004625936560341956

# Row 17
$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Hgf"
$ws.Cells.Item(17, 3).Value = "Cd44"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 0.6276006666666666
$ws.Cells.Item(17, 8).Value = 1.882802
$ws.Cells.Item(17, 9).Value = 0.01581061383652123
$ws.Cells.Item(17, 10).Value = 0.01581061383652123
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 57.695868
$ws.Cells.Item(17, 14).Value = 173.087604
$ws.Cells.Item(17, 15).Value = 0.09862711817759588
$ws.Cells.Item(17, 16).Value = 0.09862711817759588
$ws.Cells.Item(17, 17).Value = 36.20996522071199
$ws.Cells.Item(17, 18).Value = 325.889686986408
$ws.Cells.Item(17, 19).Value = 0.001559355279314912
$ws.Cells.Item(17, 20).Value = 0.001559355279314912

Write-Output "Updated rows 2-17 with ECs sending-cluster data"
